$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4663.857
$ws.Range("J29").Value = 5139.4
$ws.Range("L29").Value = 15418.2
$ws.Range("N29").Value = -15980.2
$ws.Range("H62").Value = 5990
$ws.Range("I62").Value = 3980
$ws.Range("K62").Value = 3980
$ws.Range("M62").Value = -3356
$ws.Range("H65").Value = 5990
$ws.Range("I65").Value = 3980
$ws.Range("K65").Value = 19900
$ws.Range("M65").Value = -16780
$ws.Range("H76").Value = 6711.467
$ws.Range("I76").Value = 5596.5713
$ws.Range("K76").Value = 5596.5713
$ws.Range("M76").Value = -5281.5713
$ws.Range("H79").Value = 6711.467
$ws.Range("I79").Value = 5596.5713
$ws.Range("K79").Value = 5596.5713
$ws.Range("M79").Value = -4504.5713
$ws.Range("H88").Value = 2499.3333
$ws.Range("I88").Value = 1096.5
$ws.Range("J88").Value = 3200.75
$ws.Range("K88").Value = 1096.5
$ws.Range("L88").Value = 3200.75
$ws.Range("M88").Value = -690.5
$ws.Range("N88").Value = -4012.75
$ws.Range("H91").Value = 2499.3333
$ws.Range("I91").Value = 1096.5
$ws.Range("J91").Value = 3200.75
$ws.Range("K91").Value = 1096.5
$ws.Range("L91").Value = 3200.75
$ws.Range("M91").Value = 307.5
$ws.Range("N91").Value = -6008.75
$ws.Range("H116").Value = 9999
$ws.Range("I116").Value = 9998.5
$ws.Range("K116").Value = 9998.5
$ws.Range("M116").Value = -6556.5
$ws.Range("H137").Value = 3182.1667
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3182.1667
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -14646.5001
$ws.Range("H141").Value = 4767.5
$ws.Range("I141").Value = 4645.5557
$ws.Range("J141").Value = 5133.3335
$ws.Range("K141").Value = 13936.6671
$ws.Range("L141").Value = 15400.0005
$ws.Range("M141").Value = -8756.667099999999
$ws.Range("N141").Value = -25760.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14336.46
$ws.Range("I32").Value = 5611.7896
$ws.Range("K32").Value = 5611.7896
$ws.Range("M32").Value = -5324.7896
$ws.Range("H45").Value = 4767.3335
$ws.Range("I45").Value = 2021.2
$ws.Range("K45").Value = 2021.2
$ws.Range("M45").Value = -1644.2
$ws.Range("H74").Value = 5844.923
$ws.Range("J74").Value = 7743.1113
$ws.Range("L74").Value = 7743.1113
$ws.Range("N74").Value = -9491.1113
$ws.Range("H77").Value = 5844.923
$ws.Range("J77").Value = 7743.1113
$ws.Range("L77").Value = 38715.5565
$ws.Range("N77").Value = -47451.5565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1533
$ws.Range("I20").Value = 1599.5
$ws.Range("J20").Value = 1400
$ws.Range("K20").Value = 1599.5
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = -1352.5
$ws.Range("N20").Value = -1894
$ws.Range("H86").Value = 1554.2222
$ws.Range("J86").Value = 2999.75
$ws.Range("L86").Value = 2999.75
$ws.Range("N86").Value = -5245.75
$ws.Range("H89").Value = 1554.2222
$ws.Range("J89").Value = 2999.75
$ws.Range("L89").Value = 14998.75
$ws.Range("N89").Value = -26230.75
$ws.Range("H105").Value = 3904.4
$ws.Range("I105").Value = 3222.842
$ws.Range("K105").Value = 3222.842
$ws.Range("M105").Value = -1475.842
$ws.Range("H107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 3303.5518
$ws.Range("I134").Value = 3377.9614
$ws.Range("K134").Value = 10133.8842
$ws.Range("M134").Value = -7598.8842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 605.1818
$ws.Range("I16").Value = 462
$ws.Range("K16").Value = 462
$ws.Range("M16").Value = -175
$ws.Range("H31").Value = 4397.591
$ws.Range("I31").Value = 2834.3333
$ws.Range("K31").Value = 2834.3333
$ws.Range("M31").Value = -2539.3333
$ws.Range("H34").Value = 4397.591
$ws.Range("I34").Value = 2834.3333
$ws.Range("K34").Value = 2834.3333
$ws.Range("M34").Value = -2632.3333
$ws.Range("H58").Value = 3477.2222
$ws.Range("I58").Value = 1958.4445
$ws.Range("K58").Value = 1958.4445
$ws.Range("M58").Value = -1755.4445
$ws.Range("H105").Value = 4000.818
$ws.Range("I105").Value = 3060.75
$ws.Range("J105").Value = 4538
$ws.Range("K105").Value = 3060.75
$ws.Range("L105").Value = 4538
$ws.Range("M105").Value = -1313.75
$ws.Range("N105").Value = -8032
$ws.Range("H113").Value = 605.1818
$ws.Range("I113").Value = 462
$ws.Range("K113").Value = 462
$ws.Range("M113").Value = 1708
$ws.Range("H122").Value = 8665.546
$ws.Range("I122").Value = 8687.25
$ws.Range("J122").Value = 8607.667
$ws.Range("K122").Value = 26061.75
$ws.Range("L122").Value = 25823.001
$ws.Range("M122").Value = -23611.75
$ws.Range("N122").Value = -30723.001
$ws.Range("H136").Value = 3477.2222
$ws.Range("I136").Value = 1958.4445
$ws.Range("K136").Value = 5875.333500000001
$ws.Range("M136").Value = -3325.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7714
$ws.Range("I70").Value = 7004
$ws.Range("K70").Value = 7004
$ws.Range("M70").Value = -6734
$ws.Range("H73").Value = 7714
$ws.Range("I73").Value = 7004
$ws.Range("K73").Value = 7004
$ws.Range("M73").Value = -6068
$ws.Range("H80").Value = 2942.6155
$ws.Range("I80").Value = 3258.375
$ws.Range("J80").Value = 2437.4
$ws.Range("K80").Value = 3258.375
$ws.Range("L80").Value = 2437.4
$ws.Range("M80").Value = -2260.375
$ws.Range("N80").Value = -4433.4
$ws.Range("H83").Value = 2942.6155
$ws.Range("I83").Value = 3258.375
$ws.Range("J83").Value = 2437.4
$ws.Range("K83").Value = 16291.875
$ws.Range("L83").Value = 12187
$ws.Range("M83").Value = -11299.875
$ws.Range("N83").Value = -22171
$ws.Range("H102").Value = 1560.7826
$ws.Range("I102").Value = 591.1177
$ws.Range("K102").Value = 591.1177
$ws.Range("M102").Value = 1030.8823
$ws.Range("H107").Value = 940.4783
$ws.Range("I107").Value = 822.75
$ws.Range("K107").Value = 822.75
$ws.Range("M107").Value = 1097.25
$ws.Range("H132").Value = 3502.4736
$ws.Range("I132").Value = 3103.875
$ws.Range("K132").Value = 9311.625
$ws.Range("M132").Value = -6781.625
$ws.Range("H136").Value = 25145
$ws.Range("J136").Value = 25145
$ws.Range("L136").Value = 75435
$ws.Range("N136").Value = -80535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20912
$ws.Range("H82").Value = 2277.4285
$ws.Range("I82").Value = 2178.4
$ws.Range("J82").Value = 2525
$ws.Range("K82").Value = 2178.4
$ws.Range("L82").Value = 2525
$ws.Range("M82").Value = -1817.4
$ws.Range("N82").Value = -3247
$ws.Range("H85").Value = 2277.4285
$ws.Range("I85").Value = 2178.4
$ws.Range("J85").Value = 2525
$ws.Range("K85").Value = 2178.4
$ws.Range("L85").Value = 2525
$ws.Range("M85").Value = -930.4000000000001
$ws.Range("N85").Value = -5021

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1299.2727
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 6000
$ws.Range("N100").Value = -7082
$ws.Range("H122").Value = 1286.5714
$ws.Range("I122").Value = 1286.5714
$ws.Range("K122").Value = 3859.7142
$ws.Range("M122").Value = -1409.7142
